# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the first file row (60c0026d...)
# on both the zh-cn and de-de sheets, reflecting a fresh handoff just performed.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-26 08:45:03"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-26 08:45:17"
